$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): rename the two title columns, add a new "No"
#     column in A, and relabel the "profile page" column (now column E,
#     pushing the old F header "profile definition page" to F unchanged).
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Core Profile Title"
$ws.Range("C1").Value = "Base Resource"
$ws.Range("E1").Value = "Profile Page"

# --- Column A numbering: a row was inserted for "Conformance" (row 6),
#     so every following row's running number shifts up by one. The
#     B:F content of every data row is unchanged.
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24

# Old trailing row 26 (which only held the stray "24" counter) is gone now
# that row 25 carries that value, so remove it and let the sheet shrink
# back to A1:F25.
$ws.Rows.Item(26).Delete()

# Column D (spreadsheet filename) is no longer shown.
$ws.Columns.Item(4).Hidden = $true

# Move the active selection to B2, matching the saved view state.
$ws.Range("B2").Select()
